# Update "想去人数" (interest count) figures in column F across all four
# sheets of the 北京-漫展信息 workbook, matching the regenerated data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1350
$ws.Range("F6").Value = 7705
$ws.Range("F9").Value = 2102
$ws.Range("F10").Value = 8497
$ws.Range("F13").Value = 79
$ws.Range("F14").Value = 5693
$ws.Range("F15").Value = 61
$ws.Range("F16").Value = 2650
$ws.Range("F17").Value = 1154
$ws.Range("F18").Value = 4597
$ws.Range("F23").Value = 549
$ws.Range("F24").Value = 3642
$ws.Range("F25").Value = 72
$ws.Range("F27").Value = 31
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 3136
$ws.Range("F31").Value = 217
$ws.Range("F34").Value = 336
$ws.Range("F35").Value = 916
$ws.Range("F36").Value = 673
$ws.Range("F39").Value = 2456
$ws.Range("F41").Value = 14
$ws.Range("F42").Value = 25
$ws.Range("F43").Value = 3080
$ws.Range("F45").Value = 2299
$ws.Range("F49").Value = 3

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 133
$ws.Range("F9").Value = 126

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1335

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1335
$ws.Range("F5").Value = 1350
$ws.Range("F6").Value = 7705
$ws.Range("F8").Value = 2102
$ws.Range("F9").Value = 8497
$ws.Range("F11").Value = 79
$ws.Range("F12").Value = 5693
$ws.Range("F13").Value = 61
$ws.Range("F14").Value = 2650
$ws.Range("F15").Value = 1154
$ws.Range("F16").Value = 4597
$ws.Range("F21").Value = 133
$ws.Range("F22").Value = 549
$ws.Range("F24").Value = 3642
$ws.Range("F25").Value = 72
$ws.Range("F27").Value = 31
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 3136
$ws.Range("F33").Value = 336
$ws.Range("F35").Value = 916
$ws.Range("F36").Value = 673
$ws.Range("F40").Value = 2456
$ws.Range("F42").Value = 14
$ws.Range("F43").Value = 25
$ws.Range("F44").Value = 3080
$ws.Range("F46").Value = 2300
$ws.Range("F49").Value = 126

Write-Output "Updated F-column counts across 展览, 演出, 本地生活, 全部类型 sheets."
